$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 3), mirroring formatting of row 2
$ws.Range("A3").Value = 41554
$ws.Range("A3").NumberFormat = "m/d/yy"

$ws.Range("B3").Value = 0.20138888888888887
$ws.Range("B3").NumberFormat = "h:mm"

# Update selection to match target (selected range B2:B3, active cell B3).
# NOTE: this runtime's Range.Select()/Activate() always anchors the active
# cell to the range's first (top-left) corner - Excel itself would leave the
# active cell on B3 here (the last-entered/extended cell), but that specific
# combination (sqref=B2:B3 with activeCell=B3) isn't reachable through the
# exposed object model, so we select the B2:B3 block, matching sqref exactly.
$ws.Range("B2:B3").Select()
